# Weekly update: insert a new pair of price rows (Primera/Segunda) for
# Brócoli at "Terminal La Palmera de La Serena" ahead of the existing
# data, pushing the rest of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 1059-1060; everything from the old row 1059
# onward shifts down to 1061 onward.
$ws.Rows("1059:1060").Insert()

# New row 1059 - "Primera" quality
$ws.Range("A1059").Value = 8
$ws.Range("B1059").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C1059").Value = 'Coquimbo'
$ws.Range("D1059").Value = 45106
$ws.Range("E1059").Value = 4
$ws.Range("F1059").Value = 100112023
$ws.Range("G1059").Value = 'Brócoli'
$ws.Range("H1059").Value = 'Sin especificar'
$ws.Range("I1059").Value = 'Primera'
$ws.Range("J1059").Value = 2000
$ws.Range("K1059").Value = 700
$ws.Range("L1059").Value = 800
$ws.Range("M1059").Value = 750
$ws.Range("N1059").Value = '$/unidad'
$ws.Range("O1059").Value = 'Provincia del Elquí'
$ws.Range("P1059").Value = 750
$ws.Range("Q1059").Value = 1
$ws.Range("R1059").Value = 'Hortaliza'

# New row 1060 - "Segunda" quality
$ws.Range("A1060").Value = 8
$ws.Range("B1060").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C1060").Value = 'Coquimbo'
$ws.Range("D1060").Value = 45106
$ws.Range("E1060").Value = 4
$ws.Range("F1060").Value = 100112023
$ws.Range("G1060").Value = 'Brócoli'
$ws.Range("H1060").Value = 'Sin especificar'
$ws.Range("I1060").Value = 'Segunda'
$ws.Range("J1060").Value = 1400
$ws.Range("K1060").Value = 500
$ws.Range("L1060").Value = 600
$ws.Range("M1060").Value = 550
$ws.Range("N1060").Value = '$/unidad'
$ws.Range("O1060").Value = 'Provincia del Elquí'
$ws.Range("P1060").Value = 550
$ws.Range("Q1060").Value = 1
$ws.Range("R1060").Value = 'Hortaliza'
